$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# 1) Update the F-column query timestamps on the "data" sheet.
$ws.Range("F2").Value = "2021-10-05 14:35:02.863530"
$ws.Range("F3").Value = "2021-10-05 14:35:02.863538"
$ws.Range("F4").Value = "2021-10-05 14:35:02.863541"
$ws.Range("F5").Value = "2021-10-05 14:35:02.863544"
$ws.Range("F6").Value = "2021-10-05 14:35:02.863547"
$ws.Range("F7").Value = "2021-10-05 14:35:02.863550"
$ws.Range("F8").Value = "2021-10-05 14:35:02.863552"
$ws.Range("F9").Value = "2021-10-05 14:35:02.863555"
$ws.Range("F10").Value = "2021-10-05 14:35:02.863557"
$ws.Range("F11").Value = "2021-10-05 14:35:02.863560"
$ws.Range("F12").Value = "2021-10-05 14:35:02.863563"
$ws.Range("F13").Value = "2021-10-05 14:35:02.863565"
$ws.Range("F14").Value = "2021-10-05 14:35:02.863568"
$ws.Range("F15").Value = "2021-10-05 14:35:02.863570"
$ws.Range("F16").Value = "2021-10-05 14:35:02.863573"
$ws.Range("F17").Value = "2021-10-05 14:35:02.863575"
$ws.Range("F18").Value = "2021-10-05 14:35:02.863578"
$ws.Range("F19").Value = "2021-10-05 14:35:02.863581"
$ws.Range("F20").Value = "2021-10-05 14:35:02.863583"
$ws.Range("F21").Value = "2021-10-05 14:35:02.863586"
$ws.Range("F22").Value = "2021-10-05 14:35:02.863588"
$ws.Range("F23").Value = "2021-10-05 14:35:02.863591"
$ws.Range("F24").Value = "2021-10-05 14:35:02.863593"
$ws.Range("F25").Value = "2021-10-05 14:35:02.863596"
$ws.Range("F26").Value = "2021-10-05 14:35:02.863599"
$ws.Range("F27").Value = "2021-10-05 14:35:02.863601"
$ws.Range("F28").Value = "2021-10-05 14:35:02.863604"
$ws.Range("F29").Value = "2021-10-05 14:35:02.863607"
$ws.Range("F30").Value = "2021-10-05 14:35:02.863609"
$ws.Range("F31").Value = "2021-10-05 14:35:02.863612"
$ws.Range("F32").Value = "2021-10-05 14:35:02.863614"
$ws.Range("F33").Value = "2021-10-05 14:35:02.863617"
$ws.Range("F34").Value = "2021-10-05 14:35:02.863620"
$ws.Range("F35").Value = "2021-10-05 14:35:02.863622"
$ws.Range("F36").Value = "2021-10-05 14:35:02.863625"
$ws.Range("F37").Value = "2021-10-05 14:35:02.863627"
$ws.Range("F38").Value = "2021-10-05 14:35:02.863630"

# 2) Add the new "metadata" sheet right after "data".
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

# Copy the header-row style (bold, bordered, centered) from the "data" sheet's
# header row / index column onto the metadata sheet's equivalents so they
# share the same cellXf (no new styles introduced).
$ws.Range("B1").Copy($meta.Range("B1"))
$ws.Range("B1").Copy($meta.Range("C1"))
$ws.Range("B1").Copy($meta.Range("D1"))
$ws.Range("B1").Copy($meta.Range("E1"))
$ws.Range("B1").Copy($meta.Range("F1"))
$ws.Range("B1").Copy($meta.Range("G1"))
$ws.Range("A2").Copy($meta.Range("A2"))

# Header row text.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Nephrolithiasis and Nephrocalcinosis"
$meta.Range("C2").Value = 143

# D2 ("0.56") must stay text (not coerce to a float) and keep the default
# (unstyled) cell format, so route it through a text-formula + paste-values
# round trip (using a scratch cell well outside the used range) rather than
# a direct numeric-looking Value assignment, which Excel would silently
# parse as the number 0.56.
$meta.Range("Z1").Formula = "=""0.56"""
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-09-09T06:49:55.188647Z"
$meta.Range("F2").Value = "2021-10-05 14:35:02.859892"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/143/?format=json"

$ws.Select()
